$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the two section-header rows so each cell can be addressed individually
$ws.Range("A2:E2").UnMerge()
$ws.Range("A6:E6").UnMerge()

# Reset formatting on column A: previously some label cells had center/right
# alignment styles (from the merged section headers), but in the new layout
# every column-A label cell uses the default (General) style
$ws.Range("A2:A11").Style = "Normal"

# Write the new cell contents for the whole A1:E11 table
$ws.Cells.Item(1,2).Value2 = "3 mo to 6 mo, n = 16,809"
$ws.Cells.Item(1,3).Value2 = "> 6 mo to 1 y, n = 18,894"
$ws.Cells.Item(1,4).Value2 = "> 1 y to 2 y, n =  5,421"
$ws.Cells.Item(1,5).Value2 = "> 2 y, n =  4,608"

$ws.Cells.Item(2,1).Value2 = "Overweight/obese"
$ws.Cells.Item(2,2).Value2 = ""
$ws.Cells.Item(2,3).Value2 = ""
$ws.Cells.Item(2,4).Value2 = ""
$ws.Cells.Item(2,5).Value2 = ""

$ws.Cells.Item(3,1).Value2 = "Outcome events, n"
$ws.Cells.Item(3,2).Value2 = 6103
$ws.Cells.Item(3,3).Value2 = 7114
$ws.Cells.Item(3,4).Value2 = 2249
$ws.Cells.Item(3,5).Value2 = 2386

$ws.Cells.Item(4,1).Value2 = "Total years of observation"
$ws.Cells.Item(4,2).Value2 = 35678
$ws.Cells.Item(4,3).Value2 = 39371
$ws.Cells.Item(4,4).Value2 = 10773
$ws.Cells.Item(4,5).Value2 = 9217

$ws.Cells.Item(5,1).Value2 = "Crude incidence rates per 1,000 years of observation"
$ws.Cells.Item(5,2).Value2 = 171.1
$ws.Cells.Item(5,3).Value2 = 180.7
$ws.Cells.Item(5,4).Value2 = 208.8
$ws.Cells.Item(5,5).Value2 = 258.89999999999998

$ws.Cells.Item(6,1).Value2 = "Crude incidence rate ratio (95% CI)"
$ws.Cells.Item(6,2).Value2 = "0.95 (0.89, 1.01)"
$ws.Cells.Item(6,3).Value2 = "Ref"
$ws.Cells.Item(6,4).Value2 = "1.16 (1.08, 1.24)"
$ws.Cells.Item(6,5).Value2 = "1.43 (1.34, 1.53)"

$ws.Cells.Item(7,1).Value2 = "Obese"
$ws.Cells.Item(7,2).Value2 = ""
$ws.Cells.Item(7,3).Value2 = ""
$ws.Cells.Item(7,4).Value2 = ""
$ws.Cells.Item(7,5).Value2 = ""

$ws.Cells.Item(8,1).Value2 = "Outcome events, n"
$ws.Cells.Item(8,2).Value2 = 614
$ws.Cells.Item(8,3).Value2 = 771
$ws.Cells.Item(8,4).Value2 = 265
$ws.Cells.Item(8,5).Value2 = 270

$ws.Cells.Item(9,1).Value2 = "Total years of observation"
$ws.Cells.Item(9,2).Value2 = 47922
$ws.Cells.Item(9,3).Value2 = 53848
$ws.Cells.Item(9,4).Value2 = 15415
$ws.Cells.Item(9,5).Value2 = 14654

$ws.Cells.Item(10,1).Value2 = "Crude incidence rates per 1,000 years of observation"
$ws.Cells.Item(10,2).Value2 = 12.8
$ws.Cells.Item(10,3).Value2 = 14.3
$ws.Cells.Item(10,4).Value2 = 17.2
$ws.Cells.Item(10,5).Value2 = 18.399999999999999

$ws.Cells.Item(11,1).Value2 = "Crude incidence rate ratio (95% CI)"
$ws.Cells.Item(11,2).Value2 = "0.89 (0.73, 1.09)"
$ws.Cells.Item(11,3).Value2 = "Ref"
$ws.Cells.Item(11,4).Value2 = "1.20 (0.98, 1.47)"
$ws.Cells.Item(11,5).Value2 = "1.29 (1.05, 1.57)"

# Start every data cell from a clean "Normal" style so left-over number formats
# from the previous layout (e.g. the old thousands-separator format) don't linger
$ws.Range("B1:E11").Style = "Normal"

# Apply the right-aligned (General number format) style used by most data/label cells
$ws.Range("B1:E2").HorizontalAlignment = -4152
$ws.Range("B5:E8").HorizontalAlignment = -4152
$ws.Range("B10:E11").HorizontalAlignment = -4152

# Apply the right-aligned + thousands-separator (#,##0) style used by the big counts
$ws.Range("B3:E4").NumberFormat = "#,##0"
$ws.Range("B3:E4").HorizontalAlignment = -4152
$ws.Range("B9:E9").NumberFormat = "#,##0"
$ws.Range("B9:E9").HorizontalAlignment = -4152

# Recompute "best fit" column widths for the new content
$ws.Range("A1:E11").Columns.AutoFit() | Out-Null
